# The author's commit swaps the two embedded DrawingML theme parts of
# this deck: ppt/theme/theme1.xml (used by the notes master, originally
# the stock "Office Theme" palette) and ppt/theme/theme2.xml (used by
# the slide master / all slide layouts, originally the "Integral"
# palette). After the edit theme1.xml carries the "Integral" colours
# and theme2.xml carries the "Office Theme" colours - i.e. the colour
# scheme that slides/layouts actually render with switches from
# Integral to the default Office palette. (Font scheme and format
# scheme/effect definitions are byte-identical between the two theme
# parts already, so only the 12 colour-scheme slots actually change.)
#
# Re-apply that by pushing the "Office Theme" colour values onto the
# presentation's live theme colour scheme through the object model.

$p = $ppt.ActivePresentation

$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# Index -> (slot, target "Office Theme" RGB) per ECMA-376 theme colour
# ordering (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink). RGB() isn't
# available in this host, so values are passed as the packed 0xBBGGRR
# long that PowerPoint's ColorFormat.RGB normally expects.
$colors.Colors(1).RGB  = 0x000000   # dk1      #000000
$colors.Colors(2).RGB  = 0xFFFFFF   # lt1      #FFFFFF
$colors.Colors(3).RGB  = 0x6A5444   # dk2      #44546A
$colors.Colors(4).RGB  = 0xE6E6E7   # lt2      #E7E6E6
$colors.Colors(5).RGB  = 0xD59B5B   # accent1  #5B9BD5
$colors.Colors(6).RGB  = 0x317DED   # accent2  #ED7D31
$colors.Colors(7).RGB  = 0xA5A5A5   # accent3  #A5A5A5
$colors.Colors(8).RGB  = 0x00C0FF   # accent4  #FFC000
$colors.Colors(9).RGB  = 0xC47244   # accent5  #4472C4
$colors.Colors(10).RGB = 0x47AD70   # accent6  #70AD47
$colors.Colors(11).RGB = 0xC16305   # hlink    #0563C1
$colors.Colors(12).RGB = 0x724F95   # folHlink #954F72
